# Updated transition-matrix probabilities for South Alabama_B (more games simulated).
# Values below are the refreshed cell contents per the latest simulation run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2117647058823529
$ws.Range("C2").Value = 0.5411764705882353
$ws.Range("J2").Value = 0.008823529411764706
$ws.Range("P2").Value = 0.1323529411764706
$ws.Range("S2").Value = 0.1058823529411765

# Row 3
$ws.Range("B3").Value = 0.005263157894736842
$ws.Range("C3").Value = 0.03157894736842105
$ws.Range("J3").Value = 0.005263157894736842
$ws.Range("P3").Value = 0.7210526315789474
$ws.Range("S3").Value = 0.2368421052631579

# Row 4
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.5625
$ws.Range("S4").Value = 0.375

# Row 6
$ws.Range("B6").Value = 0.0797872340425532
$ws.Range("D6").Value = 0.02659574468085106
$ws.Range("F6").Value = 0.01595744680851064
$ws.Range("J6").Value = 0.2446808510638298
$ws.Range("O6").Value = 0.005319148936170213
$ws.Range("Q6").Value = 0.148936170212766
$ws.Range("R6").Value = 0.06914893617021277
$ws.Range("S6").Value = 0.4095744680851064

# Row 7
$ws.Range("B7").Value = 0.1357466063348416
$ws.Range("D7").Value = 0.03619909502262444
$ws.Range("F7").Value = 0.03167420814479638
$ws.Range("J7").Value = 0.1493212669683258
$ws.Range("O7").Value = 0.02262443438914027
$ws.Range("Q7").Value = 0.1900452488687783
$ws.Range("R7").Value = 0.08597285067873303
$ws.Range("S7").Value = 0.3484162895927602

# Row 8
$ws.Range("B8").Value = 0.1206140350877193
$ws.Range("D8").Value = 0.02192982456140351
$ws.Range("E8").Value = 0.002192982456140351
$ws.Range("F8").Value = 0.05482456140350877
$ws.Range("J8").Value = 0.125
$ws.Range("O8").Value = 0.0131578947368421
$ws.Range("Q8").Value = 0.1885964912280702
$ws.Range("R8").Value = 0.08114035087719298
$ws.Range("S8").Value = 0.3925438596491228

# Row 9
$ws.Range("B9").Value = 0.09595959595959595
$ws.Range("D9").Value = 0.02525252525252525
$ws.Range("F9").Value = 0.09090909090909091
$ws.Range("J9").Value = 0.08080808080808081
$ws.Range("O9").Value = 0.005050505050505051
$ws.Range("Q9").Value = 0.2070707070707071
$ws.Range("R9").Value = 0.0707070707070707
$ws.Range("S9").Value = 0.4242424242424243

# Row 10
$ws.Range("B10").Value = 0.1079504011670314
$ws.Range("D10").Value = 0.01604668125455872
$ws.Range("F10").Value = 0.06345733041575492
$ws.Range("J10").Value = 0.1349380014587892
$ws.Range("O10").Value = 0.01021152443471918
$ws.Range("Q10").Value = 0.2355944566010212
$ws.Range("R10").Value = 0.07731582786287382
$ws.Range("S10").Value = 0.3544857768052516

# Row 11
$ws.Range("G11").Value = 0.1471571906354515
$ws.Range("J11").Value = 0.05351170568561873
$ws.Range("K11").Value = 0.1872909698996655
$ws.Range("L11").Value = 0.6020066889632107
$ws.Range("S11").Value = 0.01003344481605351

# Row 12
$ws.Range("G12").Value = 0.7582417582417582
$ws.Range("J12").Value = 0.2197802197802198
$ws.Range("K12").Value = 0.01098901098901099
$ws.Range("L12").Value = 0.005494505494505495
$ws.Range("S12").Value = 0.005494505494505495

# Row 13
$ws.Range("G13").Value = 0.8076923076923077
$ws.Range("J13").Value = 0.1923076923076923

# Row 14
$ws.Range("G14").Value = 0.5
$ws.Range("S14").Value = 0.5

# Row 15
$ws.Range("F15").Value = 0.01515151515151515
$ws.Range("H15").Value = 0.196969696969697
$ws.Range("I15").Value = 0.08585858585858586
$ws.Range("J15").Value = 0.3585858585858586
$ws.Range("K15").Value = 0.08585858585858586
$ws.Range("M15").Value = 0.02525252525252525
$ws.Range("O15").Value = 0.0505050505050505
$ws.Range("S15").Value = 0.1818181818181818

# Row 16
$ws.Range("F16").Value = 0.009900990099009901
$ws.Range("H16").Value = 0.1831683168316832
$ws.Range("I16").Value = 0.07425742574257425
$ws.Range("J16").Value = 0.4900990099009901
$ws.Range("K16").Value = 0.0891089108910891
$ws.Range("M16").Value = 0.02475247524752475
$ws.Range("O16").Value = 0.04950495049504951
$ws.Range("S16").Value = 0.07920792079207921

# Row 17
$ws.Range("F17").Value = 0.00576923076923077
$ws.Range("H17").Value = 0.1769230769230769
$ws.Range("I17").Value = 0.08461538461538462
$ws.Range("J17").Value = 0.4615384615384616
$ws.Range("K17").Value = 0.07115384615384615
$ws.Range("M17").Value = 0.025
$ws.Range("N17").Value = 0.001923076923076923
$ws.Range("O17").Value = 0.0673076923076923
$ws.Range("S17").Value = 0.1057692307692308

# Row 18
$ws.Range("F18").Value = 0.01612903225806452
$ws.Range("H18").Value = 0.1720430107526882
$ws.Range("I18").Value = 0.1075268817204301
$ws.Range("J18").Value = 0.4569892473118279
$ws.Range("K18").Value = 0.08602150537634409
$ws.Range("M18").Value = 0.02150537634408602
$ws.Range("O18").Value = 0.05376344086021505
$ws.Range("S18").Value = 0.08602150537634409

# Row 19
$ws.Range("F19").Value = 0.01383238405207486
$ws.Range("H19").Value = 0.2148087876322213
$ws.Range("I19").Value = 0.0821806346623271
$ws.Range("J19").Value = 0.3938161106590724
$ws.Range("K19").Value = 0.1220504475183076
$ws.Range("M19").Value = 0.02115541090317331
$ws.Range("O19").Value = 0.06102522375915378
$ws.Range("S19").Value = 0.09113100081366965

